$d = $word.ActiveDocument
$count = 0

if ($d.Content.Find.Execute("2025-06-12 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-06-13 Friday", 2)) { $count++ }
if ($d.Content.Find.Execute("9+66=75", $true, $false, $false, $false, $false, $true, 1, $false, "80-47=33", 2)) { $count++ }
if ($d.Content.Find.Execute("8+11=19", $true, $false, $false, $false, $false, $true, 1, $false, "10+52=62", 2)) { $count++ }
if ($d.Content.Find.Execute("40-30=10", $true, $false, $false, $false, $false, $true, 1, $false, "44+43=87", 2)) { $count++ }
if ($d.Content.Find.Execute("1+10=11", $true, $false, $false, $false, $false, $true, 1, $false, "84-7=77", 2)) { $count++ }
if ($d.Content.Find.Execute("60-14=46", $true, $false, $false, $false, $false, $true, 1, $false, "9+2=11", 2)) { $count++ }
if ($d.Content.Find.Execute("25+11=36", $true, $false, $false, $false, $false, $true, 1, $false, "15+8=23", 2)) { $count++ }
if ($d.Content.Find.Execute("55+17=72", $true, $false, $false, $false, $false, $true, 1, $false, "77-23=54", 2)) { $count++ }
if ($d.Content.Find.Execute("73-34=39", $true, $false, $false, $false, $false, $true, 1, $false, "41+43=84", 2)) { $count++ }
if ($d.Content.Find.Execute("84-1=83", $true, $false, $false, $false, $false, $true, 1, $false, "71-54=17", 2)) { $count++ }
if ($d.Content.Find.Execute("98-90=8", $true, $false, $false, $false, $false, $true, 1, $false, "42+7=49", 2)) { $count++ }
if ($d.Content.Find.Execute("11+79=90", $true, $false, $false, $false, $false, $true, 1, $false, "91-55=36", 2)) { $count++ }
if ($d.Content.Find.Execute("43+35=78", $true, $false, $false, $false, $false, $true, 1, $false, "38+24=62", 2)) { $count++ }
if ($d.Content.Find.Execute("36-7=29", $true, $false, $false, $false, $false, $true, 1, $false, "72-25=47", 2)) { $count++ }
if ($d.Content.Find.Execute("11+8=19", $true, $false, $false, $false, $false, $true, 1, $false, "5+22=27", 2)) { $count++ }
if ($d.Content.Find.Execute("19+5=24", $true, $false, $false, $false, $false, $true, 1, $false, "95-91=4", 2)) { $count++ }
if ($d.Content.Find.Execute("82-9=73", $true, $false, $false, $false, $false, $true, 1, $false, "73-59=14", 2)) { $count++ }
if ($d.Content.Find.Execute("14+39=53", $true, $false, $false, $false, $false, $true, 1, $false, "80-0=80", 2)) { $count++ }
if ($d.Content.Find.Execute("52+9=61", $true, $false, $false, $false, $false, $true, 1, $false, "93+1=94", 2)) { $count++ }
if ($d.Content.Find.Execute("91-50=41", $true, $false, $false, $false, $false, $true, 1, $false, "9+37=46", 2)) { $count++ }
if ($d.Content.Find.Execute("37+21=58", $true, $false, $false, $false, $false, $true, 1, $false, "44+37=81", 2)) { $count++ }
if ($d.Content.Find.Execute("98-66=32", $true, $false, $false, $false, $false, $true, 1, $false, "10+16=26", 2)) { $count++ }
if ($d.Content.Find.Execute("24-16=8", $true, $false, $false, $false, $false, $true, 1, $false, "54-7=47", 2)) { $count++ }
if ($d.Content.Find.Execute("55-19=36", $true, $false, $false, $false, $false, $true, 1, $false, "58-14=44", 2)) { $count++ }
if ($d.Content.Find.Execute("71-32=39", $true, $false, $false, $false, $false, $true, 1, $false, "77-34=43", 2)) { $count++ }
if ($d.Content.Find.Execute("49-24=25", $true, $false, $false, $false, $false, $true, 1, $false, "22+48=70", 2)) { $count++ }
if ($d.Content.Find.Execute("83-9=74", $true, $false, $false, $false, $false, $true, 1, $false, "67-8=59", 2)) { $count++ }
if ($d.Content.Find.Execute("84-78=6", $true, $false, $false, $false, $false, $true, 1, $false, "72-1=71", 2)) { $count++ }
if ($d.Content.Find.Execute("65-52=13", $true, $false, $false, $false, $false, $true, 1, $false, "41+54=95", 2)) { $count++ }
if ($d.Content.Find.Execute("20+58=78", $true, $false, $false, $false, $false, $true, 1, $false, "27+15=42", 2)) { $count++ }
if ($d.Content.Find.Execute("34+24=58", $true, $false, $false, $false, $false, $true, 1, $false, "72+16=88", 2)) { $count++ }
if ($d.Content.Find.Execute("68-39=29", $true, $false, $false, $false, $false, $true, 1, $false, "26+4=30", 2)) { $count++ }
if ($d.Content.Find.Execute("16-14=2", $true, $false, $false, $false, $false, $true, 1, $false, "35+34=69", 2)) { $count++ }
if ($d.Content.Find.Execute("37+35=72", $true, $false, $false, $false, $false, $true, 1, $false, "21+32=53", 2)) { $count++ }
if ($d.Content.Find.Execute("46-26=20", $true, $false, $false, $false, $false, $true, 1, $false, "89-67=22", 2)) { $count++ }
if ($d.Content.Find.Execute("80-35=45", $true, $false, $false, $false, $false, $true, 1, $false, "45+12=57", 2)) { $count++ }
if ($d.Content.Find.Execute("5+32=37", $true, $false, $false, $false, $false, $true, 1, $false, "66+27=93", 2)) { $count++ }
if ($d.Content.Find.Execute("40+5=45", $true, $false, $false, $false, $false, $true, 1, $false, "75+23=98", 2)) { $count++ }
if ($d.Content.Find.Execute("15+24=39", $true, $false, $false, $false, $false, $true, 1, $false, "23+53=76", 2)) { $count++ }
if ($d.Content.Find.Execute("61-12=49", $true, $false, $false, $false, $false, $true, 1, $false, "79+6=85", 2)) { $count++ }
if ($d.Content.Find.Execute("78-3=75", $true, $false, $false, $false, $false, $true, 1, $false, "15+79=94", 2)) { $count++ }
if ($d.Content.Find.Execute("15+34=49", $true, $false, $false, $false, $false, $true, 1, $false, "74-62=12", 2)) { $count++ }
if ($d.Content.Find.Execute("96+3=99", $true, $false, $false, $false, $false, $true, 1, $false, "12+61=73", 2)) { $count++ }
if ($d.Content.Find.Execute("89-44=45", $true, $false, $false, $false, $false, $true, 1, $false, "13+66=79", 2)) { $count++ }
if ($d.Content.Find.Execute("27+29=56", $true, $false, $false, $false, $false, $true, 1, $false, "35+37=72", 2)) { $count++ }
if ($d.Content.Find.Execute("53+40=93", $true, $false, $false, $false, $false, $true, 1, $false, "51-33=18", 2)) { $count++ }
if ($d.Content.Find.Execute("40+38=78", $true, $false, $false, $false, $false, $true, 1, $false, "43+54=97", 2)) { $count++ }
if ($d.Content.Find.Execute("45+4=49", $true, $false, $false, $false, $false, $true, 1, $false, "1+96=97", 2)) { $count++ }
if ($d.Content.Find.Execute("58-45=13", $true, $false, $false, $false, $false, $true, 1, $false, "23-17=6", 2)) { $count++ }
if ($d.Content.Find.Execute("39+33=72", $true, $false, $false, $false, $false, $true, 1, $false, "25+58=83", 2)) { $count++ }
if ($d.Content.Find.Execute("7+84=91", $true, $false, $false, $false, $false, $true, 1, $false, "80+8=88", 2)) { $count++ }
if ($d.Content.Find.Execute("19-13=6", $true, $false, $false, $false, $false, $true, 1, $false, "10-7=3", 2)) { $count++ }
if ($d.Content.Find.Execute("50+38=88", $true, $false, $false, $false, $false, $true, 1, $false, "88-67=21", 2)) { $count++ }
if ($d.Content.Find.Execute("83-13=70", $true, $false, $false, $false, $false, $true, 1, $false, "31+45=76", 2)) { $count++ }
if ($d.Content.Find.Execute("48+34=82", $true, $false, $false, $false, $false, $true, 1, $false, "52+33=85", 2)) { $count++ }
if ($d.Content.Find.Execute("2+34=36", $true, $false, $false, $false, $false, $true, 1, $false, "37-23=14", 2)) { $count++ }
if ($d.Content.Find.Execute("75+7=82", $true, $false, $false, $false, $false, $true, 1, $false, "92-78=14", 2)) { $count++ }
if ($d.Content.Find.Execute("31-27=4", $true, $false, $false, $false, $false, $true, 1, $false, "41+54=95", 2)) { $count++ }
if ($d.Content.Find.Execute("95-38=57", $true, $false, $false, $false, $false, $true, 1, $false, "0+39=39", 2)) { $count++ }
if ($d.Content.Find.Execute("44-25=19", $true, $false, $false, $false, $false, $true, 1, $false, "33+23=56", 2)) { $count++ }
if ($d.Content.Find.Execute("53-38=15", $true, $false, $false, $false, $false, $true, 1, $false, "41-3=38", 2)) { $count++ }
if ($d.Content.Find.Execute("42+34=76", $true, $false, $false, $false, $false, $true, 1, $false, "73-18=55", 2)) { $count++ }
if ($d.Content.Find.Execute("84-57=27", $true, $false, $false, $false, $false, $true, 1, $false, "57-38=19", 2)) { $count++ }
if ($d.Content.Find.Execute("21+5=26", $true, $false, $false, $false, $false, $true, 1, $false, "83-73=10", 2)) { $count++ }
if ($d.Content.Find.Execute("5+89=94", $true, $false, $false, $false, $false, $true, 1, $false, "29+67=96", 2)) { $count++ }
if ($d.Content.Find.Execute("94+0=94", $true, $false, $false, $false, $false, $true, 1, $false, "50-38=12", 2)) { $count++ }
if ($d.Content.Find.Execute("45+1=46", $true, $false, $false, $false, $false, $true, 1, $false, "37+42=79", 2)) { $count++ }
if ($d.Content.Find.Execute("11+45=56", $true, $false, $false, $false, $false, $true, 1, $false, "67-39=28", 2)) { $count++ }
if ($d.Content.Find.Execute("66+16=82", $true, $false, $false, $false, $false, $true, 1, $false, "49+7=56", 2)) { $count++ }
if ($d.Content.Find.Execute("51+42=93", $true, $false, $false, $false, $false, $true, 1, $false, "75-47=28", 2)) { $count++ }
if ($d.Content.Find.Execute("51-47=4", $true, $false, $false, $false, $false, $true, 1, $false, "50-32=18", 2)) { $count++ }
if ($d.Content.Find.Execute("24+10=34", $true, $false, $false, $false, $false, $true, 1, $false, "32-24=8", 2)) { $count++ }
if ($d.Content.Find.Execute("73-51=22", $true, $false, $false, $false, $false, $true, 1, $false, "49-49=0", 2)) { $count++ }
if ($d.Content.Find.Execute("43+46=89", $true, $false, $false, $false, $false, $true, 1, $false, "21+26=47", 2)) { $count++ }
if ($d.Content.Find.Execute("18+26=44", $true, $false, $false, $false, $false, $true, 1, $false, "37-3=34", 2)) { $count++ }
if ($d.Content.Find.Execute("49-29=20", $true, $false, $false, $false, $false, $true, 1, $false, "8-6=2", 2)) { $count++ }
if ($d.Content.Find.Execute("73+16=89", $true, $false, $false, $false, $false, $true, 1, $false, "92-20=72", 2)) { $count++ }
if ($d.Content.Find.Execute("46-44=2", $true, $false, $false, $false, $false, $true, 1, $false, "8+75=83", 2)) { $count++ }
if ($d.Content.Find.Execute("77-73=4", $true, $false, $false, $false, $false, $true, 1, $false, "20+8=28", 2)) { $count++ }
if ($d.Content.Find.Execute("23+35=58", $true, $false, $false, $false, $false, $true, 1, $false, "18+60=78", 2)) { $count++ }
if ($d.Content.Find.Execute("24+64=88", $true, $false, $false, $false, $false, $true, 1, $false, "89-21=68", 2)) { $count++ }
if ($d.Content.Find.Execute("76-43=33", $true, $false, $false, $false, $false, $true, 1, $false, "63-27=36", 2)) { $count++ }
if ($d.Content.Find.Execute("40-24=16", $true, $false, $false, $false, $false, $true, 1, $false, "11+59=70", 2)) { $count++ }
if ($d.Content.Find.Execute("76-30=46", $true, $false, $false, $false, $false, $true, 1, $false, "18+41=59", 2)) { $count++ }
if ($d.Content.Find.Execute("90-16=74", $true, $false, $false, $false, $false, $true, 1, $false, "36+52=88", 2)) { $count++ }
if ($d.Content.Find.Execute("71-68=3", $true, $false, $false, $false, $false, $true, 1, $false, "56+6=62", 2)) { $count++ }
if ($d.Content.Find.Execute("58-10=48", $true, $false, $false, $false, $false, $true, 1, $false, "73-24=49", 2)) { $count++ }
if ($d.Content.Find.Execute("58-12=46", $true, $false, $false, $false, $false, $true, 1, $false, "59+40=99", 2)) { $count++ }
if ($d.Content.Find.Execute("60+33=93", $true, $false, $false, $false, $false, $true, 1, $false, "2+2=4", 2)) { $count++ }
if ($d.Content.Find.Execute("27+49=76", $true, $false, $false, $false, $false, $true, 1, $false, "85+10=95", 2)) { $count++ }
if ($d.Content.Find.Execute("57-42=15", $true, $false, $false, $false, $false, $true, 1, $false, "84-40=44", 2)) { $count++ }
if ($d.Content.Find.Execute("12-5=7", $true, $false, $false, $false, $false, $true, 1, $false, "93-15=78", 2)) { $count++ }
if ($d.Content.Find.Execute("71+18=89", $true, $false, $false, $false, $false, $true, 1, $false, "53+35=88", 2)) { $count++ }
if ($d.Content.Find.Execute("10+12=22", $true, $false, $false, $false, $false, $true, 1, $false, "72-18=54", 2)) { $count++ }
if ($d.Content.Find.Execute("84+13=97", $true, $false, $false, $false, $false, $true, 1, $false, "60-22=38", 2)) { $count++ }
if ($d.Content.Find.Execute("86-6=80", $true, $false, $false, $false, $false, $true, 1, $false, "26+11=37", 2)) { $count++ }
if ($d.Content.Find.Execute("86-85=1", $true, $false, $false, $false, $false, $true, 1, $false, "1+31=32", 2)) { $count++ }
if ($d.Content.Find.Execute("0+30=30", $true, $false, $false, $false, $false, $true, 1, $false, "14-2=12", 2)) { $count++ }
if ($d.Content.Find.Execute("51-40=11", $true, $false, $false, $false, $false, $true, 1, $false, "35-11=24", 2)) { $count++ }
if ($d.Content.Find.Execute("75-23=52", $true, $false, $false, $false, $false, $true, 1, $false, "2+9=11", 2)) { $count++ }
if ($d.Content.Find.Execute("85-34=51", $true, $false, $false, $false, $false, $true, 1, $false, "36-16=20", 2)) { $count++ }

Write-Output "Replacements applied: $count / 101"
